$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 4 (ALC)
$ws.Cells.Item(4, 8).Value = 51.5  # H4: was 51.2
$ws.Cells.Item(4, 9).Value = 45  # I4: was 47.5
$ws.Cells.Item(4, 11).Value = 45  # K4: was 47.5
$ws.Cells.Item(4, 13).Value = 69  # M4: was 66.5

# row 9 (ALC)
$ws.Cells.Item(9, 8).Value = 83.5  # H9: was 83.09999999999999
$ws.Cells.Item(9, 9).Value = 74  # I9: was 71.71429000000001
$ws.Cells.Item(9, 10).Value = 150  # J9: was 109.666664
$ws.Cells.Item(9, 11).Value = 74  # K9: was 71.71429000000001
$ws.Cells.Item(9, 12).Value = 150  # L9: was 109.666664
$ws.Cells.Item(9, 13).Value = 95  # M9: was 97.28570999999999
$ws.Cells.Item(9, 14).Value = -488  # N9: was -447.666664

# row 19 (ALC)
$ws.Cells.Item(19, 8).Value = 1557.2858  # H19: was 1500.1428
$ws.Cells.Item(19, 9).Value = 1637.75  # I19: was 1537.75
$ws.Cells.Item(19, 11).Value = 1637.75  # K19: was 1537.75
$ws.Cells.Item(19, 13).Value = -1462.75  # M19: was -1362.75

# row 43 (ALC)
$ws.Cells.Item(43, 8).Value = 2313.25  # H43: was 1994.6
$ws.Cells.Item(43, 9).Value = 1917.6666  # I43: was 1618.25
$ws.Cells.Item(43, 11).Value = 1917.6666  # K43: was 1618.25
$ws.Cells.Item(43, 13).Value = -1848.6666  # M43: was -1549.25

# row 55 (ALC)
$ws.Cells.Item(55, 8).Value = 726.1111  # H55: was 467.1
$ws.Cells.Item(55, 9).Value = 592.3333  # I55: was 566.75
$ws.Cells.Item(55, 10).Value = 793  # J55: was 400.66666
$ws.Cells.Item(55, 11).Value = 592.3333  # K55: was 566.75
$ws.Cells.Item(55, 12).Value = 793  # L55: was 400.66666
$ws.Cells.Item(55, 13).Value = -378.3333  # M55: was -352.75
$ws.Cells.Item(55, 14).Value = -1221  # N55: was -828.66666

# row 132 (ALC)
$ws.Cells.Item(132, 8).Value = 1306.25  # H132: was 1927.7778
$ws.Cells.Item(132, 9).Value = 1072.3572  # I132: was 1054.2
$ws.Cells.Item(132, 10).Value = 2943.5  # J132: was 6295.6665
$ws.Cells.Item(132, 11).Value = 3217.0716  # K132: was 3162.6
$ws.Cells.Item(132, 12).Value = 8830.5  # L132: was 18886.9995
$ws.Cells.Item(132, 13).Value = -687.0715999999998  # M132: was -632.6000000000004
$ws.Cells.Item(132, 14).Value = -13890.5  # N132: was -23946.9995

# row 137 (ALC)
$ws.Cells.Item(137, 8).Value = 2154.4  # H137: was 2154.7
$ws.Cells.Item(137, 9).Value = 728.8  # I137: was 729.4
$ws.Cells.Item(137, 11).Value = 2186.4  # K137: was 2188.2
$ws.Cells.Item(137, 13).Value = 363.6000000000004  # M137: was 361.8000000000002

$ws = $wb.Worksheets.Item("ARM")
# row 122 (ARM)
$ws.Cells.Item(122, 8).Value = 2449.8  # H122: was 2537.25
$ws.Cells.Item(122, 9).Value = 2449.8  # I122: was 2537.25
$ws.Cells.Item(122, 11).Value = 7349.400000000001  # K122: was 7611.75
$ws.Cells.Item(122, 13).Value = -4899.400000000001  # M122: was -5161.75

$ws = $wb.Worksheets.Item("BSM")
# row 7 (BSM)
$ws.Cells.Item(7, 8).Value = 8337055.5  # H7: was 8827441
$ws.Cells.Item(7, 10).Value = 8000399.5  # J7: was 10000375
$ws.Cells.Item(7, 12).Value = 8000399.5  # L7: was 10000375
$ws.Cells.Item(7, 14).Value = -8000625.5  # N7: was -10000601

# row 80 (BSM)
$ws.Cells.Item(80, 8).Value = 557.73334  # H80: was 542
$ws.Cells.Item(80, 9).Value = 274  # I80: was 275
$ws.Cells.Item(80, 10).Value = 1125.2  # J80: was 1209.5
$ws.Cells.Item(80, 11).Value = 274  # K80: was 275
$ws.Cells.Item(80, 12).Value = 1125.2  # L80: was 1209.5
$ws.Cells.Item(80, 13).Value = 724  # M80: was 723
$ws.Cells.Item(80, 14).Value = -3121.2  # N80: was -3205.5

# row 83 (BSM)
$ws.Cells.Item(83, 8).Value = 557.73334  # H83: was 542
$ws.Cells.Item(83, 9).Value = 274  # I83: was 275
$ws.Cells.Item(83, 10).Value = 1125.2  # J83: was 1209.5
$ws.Cells.Item(83, 11).Value = 1370  # K83: was 1375
$ws.Cells.Item(83, 12).Value = 5626  # L83: was 6047.5
$ws.Cells.Item(83, 13).Value = 3622  # M83: was 3617
$ws.Cells.Item(83, 14).Value = -15610  # N83: was -16031.5

# row 99 (BSM)
$ws.Cells.Item(99, 8).Value = 2490.6  # H99: was 2743.25
$ws.Cells.Item(99, 9).Value = 1485  # I99: was 1487.5
$ws.Cells.Item(99, 11).Value = 1485  # K99: was 1487.5
$ws.Cells.Item(99, 13).Value = 13  # M99: was 10.5

$ws = $wb.Worksheets.Item("CRP")
# row 41 (CRP)
$ws.Cells.Item(41, 8).Value = 42908.125  # H41: was 34676.39
$ws.Cells.Item(41, 10).Value = 45466.43  # J41: was 35116.227
$ws.Cells.Item(41, 12).Value = 45466.43  # L41: was 35116.227
$ws.Cells.Item(41, 14).Value = -46322.43  # N41: was -35972.227

# row 122 (CRP)
$ws.Cells.Item(122, 8).Value = 2750.1428  # H122: was 2886
$ws.Cells.Item(122, 9).Value = 2447.3333  # I122: was 2886
$ws.Cells.Item(122, 10).Value = 4567  # J122: was 0
$ws.Cells.Item(122, 11).Value = 7341.999899999999  # K122: was 8658
$ws.Cells.Item(122, 12).Value = 13701  # L122: was 0
$ws.Cells.Item(122, 13).Value = -4891.999899999999  # M122: was -6208
$ws.Cells.Item(122, 14).Value = -18601  # N122: was None

$ws = $wb.Worksheets.Item("CUL")
# row 15 (CUL)
$ws.Cells.Item(15, 8).Value = 159.8  # H15: was 399.83334
$ws.Cells.Item(15, 9).Value = 99  # I15: was 799.5
$ws.Cells.Item(15, 10).Value = 175  # J15: was 200
$ws.Cells.Item(15, 11).Value = 297  # K15: was 2398.5
$ws.Cells.Item(15, 12).Value = 525  # L15: was 600
$ws.Cells.Item(15, 13).Value = -157  # M15: was -2258.5
$ws.Cells.Item(15, 14).Value = -805  # N15: was -880

# row 55 (CUL)
$ws.Cells.Item(55, 8).Value = 19333  # H55: was 7777.6665
$ws.Cells.Item(55, 10).Value = 19333  # J55: was 7777.6665
$ws.Cells.Item(55, 12).Value = 57999  # L55: was 23332.9995
$ws.Cells.Item(55, 14).Value = -58353  # N55: was -23686.9995

# row 68 (CUL)
$ws.Cells.Item(68, 8).Value = 3216.75  # H68: was 3922.3333
$ws.Cells.Item(68, 9).Value = 1500  # I68: was 0
$ws.Cells.Item(68, 10).Value = 3789  # J68: was 3922.3333
$ws.Cells.Item(68, 11).Value = 4500  # K68: was 0
$ws.Cells.Item(68, 12).Value = 11367  # L68: was 11766.9999
$ws.Cells.Item(68, 13).Value = -3689  # M68: was None
$ws.Cells.Item(68, 14).Value = -12989  # N68: was -13388.9999

# row 71 (CUL)
$ws.Cells.Item(71, 8).Value = 3216.75  # H71: was 3922.3333
$ws.Cells.Item(71, 9).Value = 1500  # I71: was 0
$ws.Cells.Item(71, 10).Value = 3789  # J71: was 3922.3333
$ws.Cells.Item(71, 11).Value = 13500  # K71: was 0
$ws.Cells.Item(71, 12).Value = 34101  # L71: was 35300.9997
$ws.Cells.Item(71, 13).Value = -9444  # M71: was None
$ws.Cells.Item(71, 14).Value = -42213  # N71: was -43412.9997

# row 125 (CUL)
$ws.Cells.Item(125, 8).Value = 7614.75  # H125: was 7493.3335
$ws.Cells.Item(125, 9).Value = 7249.5  # I125: was 7250
$ws.Cells.Item(125, 11).Value = 21748.5  # K125: was 21750
$ws.Cells.Item(125, 13).Value = -16828.5  # M125: was -16830

# row 132 (CUL)
$ws.Cells.Item(132, 8).Value = 3249.4  # H132: was 2513.9285
$ws.Cells.Item(132, 9).Value = 3062.25  # I132: was 2400.182
$ws.Cells.Item(132, 10).Value = 3998  # J132: was 2931
$ws.Cells.Item(132, 11).Value = 27560.25  # K132: was 21601.638
$ws.Cells.Item(132, 12).Value = 35982  # L132: was 26379
$ws.Cells.Item(132, 13).Value = -25030.25  # M132: was -19071.638
$ws.Cells.Item(132, 14).Value = -41042  # N132: was -31439

$ws = $wb.Worksheets.Item("GSM")
# row 2 (GSM)
$ws.Cells.Item(2, 8).Value = 260  # H2: was 259.5
$ws.Cells.Item(2, 9).Value = 235.55556  # I2: was 223
$ws.Cells.Item(2, 10).Value = 284.44446  # J2: was 305.125
$ws.Cells.Item(2, 11).Value = 235.55556  # K2: was 223
$ws.Cells.Item(2, 12).Value = 284.44446  # L2: was 305.125
$ws.Cells.Item(2, 13).Value = -122.55556  # M2: was -110
$ws.Cells.Item(2, 14).Value = -510.44446  # N2: was -531.125

# row 25 (GSM)
$ws.Cells.Item(25, 8).Value = 0  # H25: was 1000
$ws.Cells.Item(25, 10).Value = 0  # J25: was 1000
$ws.Cells.Item(25, 12).ClearContents()  # L25: was 1000
$ws.Cells.Item(25, 14).Value = 0  # N25: was -2058

# row 26 (GSM)
$ws.Cells.Item(26, 8).Value = 10038  # H26: was 0
$ws.Cells.Item(26, 9).Value = 10038  # I26: was 0
$ws.Cells.Item(26, 11).Value = 10038  # K26: was 0
$ws.Cells.Item(26, 13).Value = -9758  # M26: was None

# row 50 (GSM)
$ws.Cells.Item(50, 8).Value = 10038  # H50: was 0
$ws.Cells.Item(50, 9).Value = 10038  # I50: was 0
$ws.Cells.Item(50, 11).Value = 10038  # K50: was 0
$ws.Cells.Item(50, 13).Value = -9540  # M50: was None

# row 105 (GSM)
$ws.Cells.Item(105, 8).Value = 35000  # H105: was 22500
$ws.Cells.Item(105, 10).Value = 35000  # J105: was 22500
$ws.Cells.Item(105, 12).Value = 35000  # L105: was 22500
$ws.Cells.Item(105, 14).Value = -41988  # N105: was -29488

# row 107 (GSM)
$ws.Cells.Item(107, 8).Value = 0  # H107: was 214.66667
$ws.Cells.Item(107, 9).Value = 0  # I107: was 222
$ws.Cells.Item(107, 10).Value = 0  # J107: was 200
$ws.Cells.Item(107, 11).Value = 0  # K107: was 222
$ws.Cells.Item(107, 12).ClearContents()  # L107: was 200
$ws.Cells.Item(107, 13).ClearContents()  # M107: was 1698
$ws.Cells.Item(107, 14).Value = 0  # N107: was -4040

$ws = $wb.Worksheets.Item("LTW")
# row 2 (LTW)
$ws.Cells.Item(2, 8).Value = 11562.5  # H2: was 11999.5
$ws.Cells.Item(2, 10).Value = 11562.5  # J2: was 11999.5
$ws.Cells.Item(2, 12).Value = 11562.5  # L2: was 11999.5
$ws.Cells.Item(2, 14).Value = -11786.5  # N2: was -12223.5

# row 46 (LTW)
$ws.Cells.Item(46, 8).Value = 6714.2856  # H46: was 6055.5625
$ws.Cells.Item(46, 9).Value = 4400  # I46: was 3998.3333
$ws.Cells.Item(46, 10).Value = 8000  # J46: was 7289.9
$ws.Cells.Item(46, 11).Value = 4400  # K46: was 3998.3333
$ws.Cells.Item(46, 12).Value = 8000  # L46: was 7289.9
$ws.Cells.Item(46, 13).Value = -4212  # M46: was -3810.3333
$ws.Cells.Item(46, 14).Value = -8376  # N46: was -7665.9

# row 132 (LTW)
$ws.Cells.Item(132, 8).Value = 1749.75  # H132: was 1950
$ws.Cells.Item(132, 9).Value = 1749.75  # I132: was 1950
$ws.Cells.Item(132, 11).Value = 5249.25  # K132: was 5850
$ws.Cells.Item(132, 13).Value = -2719.25  # M132: was -3320

$ws = $wb.Worksheets.Item("WVR")
# row 4 (WVR)
$ws.Cells.Item(4, 8).Value = 8452.916999999999  # H4: was 9175.909
$ws.Cells.Item(4, 9).Value = 9212.272000000001  # I4: was 10083.5
$ws.Cells.Item(4, 11).Value = 9212.272000000001  # K4: was 10083.5
$ws.Cells.Item(4, 13).Value = -9099.272000000001  # M4: was -9970.5

# row 6 (WVR)
$ws.Cells.Item(6, 8).Value = 6668583.5  # H6: was 10000375
$ws.Cells.Item(6, 10).Value = 10002500  # J6: was 20000000
$ws.Cells.Item(6, 12).Value = 10002500  # L6: was 20000000
$ws.Cells.Item(6, 14).Value = -10002730  # N6: was -20000230

# row 55 (WVR)
$ws.Cells.Item(55, 8).Value = 580.5  # H55: was 602.6667
$ws.Cells.Item(55, 10).Value = 514.5  # J55: was 515
$ws.Cells.Item(55, 12).Value = 514.5  # L55: was 515
$ws.Cells.Item(55, 14).Value = -1068.5  # N55: was -1069

# row 96 (WVR)
$ws.Cells.Item(96, 8).Value = 1600.5  # H96: was 1834.3334
$ws.Cells.Item(96, 9).Value = 1634  # I96: was 2001.5
$ws.Cells.Item(96, 11).Value = 1634  # K96: was 2001.5
$ws.Cells.Item(96, 13).Value = -261  # M96: was -628.5

# row 132 (WVR)
$ws.Cells.Item(132, 8).Value = 2025.7273  # H132: was 2209.3333
$ws.Cells.Item(132, 9).Value = 2064.7778  # I132: was 2312
$ws.Cells.Item(132, 11).Value = 6194.3334  # K132: was 6936
$ws.Cells.Item(132, 13).Value = -3664.3334  # M132: was -4406
